# chore: update Sheets via scheduled runner
# Refresh cached market-price / profit figures (columns H-N) across the
# per-job Leve tables. Values below are the updated prices/profits pulled
# by the scheduled runner; row/column layout is unchanged.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 128.3
$ws.Range("I12").Value = 128.3
$ws.Range("K12").Value = 128.3
$ws.Range("M12").Value = 41.69999999999999

$ws.Range("H40").Value = 6666.6665
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H62").Value = 7688.7896
$ws.Range("I62").Value = 7339.8
$ws.Range("K62").Value = 7339.8
$ws.Range("M62").Value = -6715.8

$ws.Range("H65").Value = 7688.7896
$ws.Range("I65").Value = 7339.8
$ws.Range("K65").Value = 36699
$ws.Range("M65").Value = -33579

$ws.Range("H100").Value = 1003.5
$ws.Range("I100").Value = 874.1429000000001
$ws.Range("J100").Value = 1184.6
$ws.Range("K100").Value = 874.1429000000001
$ws.Range("L100").Value = 1184.6
$ws.Range("M100").Value = -333.1429000000001
$ws.Range("N100").Value = -2266.6

$ws.Range("H111").Value = 2695
$ws.Range("I111").Value = 2029
$ws.Range("J111").Value = 4027
$ws.Range("K111").Value = 6087
$ws.Range("L111").Value = 12081
$ws.Range("M111").Value = -3020
$ws.Range("N111").Value = -18215

$ws.Range("H115").Value = 1064.75
$ws.Range("I115").Value = 1064.75
$ws.Range("K115").Value = 3194.25
$ws.Range("M115").Value = -1627.25

$ws.Range("H129").Value = 1731.28
$ws.Range("I129").Value = 1066.8462
$ws.Range("J129").Value = 2451.0833
$ws.Range("K129").Value = 3200.5386
$ws.Range("L129").Value = 7353.249899999999
$ws.Range("M129").Value = 1799.4614
$ws.Range("N129").Value = -17353.2499

$ws.Range("H138").Value = 3842.4
$ws.Range("J138").Value = 4903.759
$ws.Range("L138").Value = 14711.277
$ws.Range("N138").Value = -24991.277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 25012500
$ws.Range("J9").Value = 25000
$ws.Range("L9").Value = 25000
$ws.Range("N9").Value = -25340

$ws.Range("H20").Value = 25012500
$ws.Range("J20").Value = 25000
$ws.Range("L20").Value = 25000
$ws.Range("N20").Value = -25540

$ws.Range("H45").Value = 1633.5416
$ws.Range("I45").Value = 1477.1
$ws.Range("K45").Value = 1477.1
$ws.Range("M45").Value = -1100.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 572.1111
$ws.Range("I22").Value = 326.16666
$ws.Range("K22").Value = 326.16666
$ws.Range("M22").Value = -153.16666

$ws.Range("H86").Value = 3428.158
$ws.Range("I86").Value = 2106.1428
$ws.Range("J86").Value = 7129.8
$ws.Range("K86").Value = 2106.1428
$ws.Range("L86").Value = 7129.8
$ws.Range("M86").Value = -983.1428000000001
$ws.Range("N86").Value = -9375.799999999999

$ws.Range("H89").Value = 3428.158
$ws.Range("I89").Value = 2106.1428
$ws.Range("J89").Value = 7129.8
$ws.Range("K89").Value = 10530.714
$ws.Range("L89").Value = 35649
$ws.Range("M89").Value = -4914.714
$ws.Range("N89").Value = -46881

$ws.Range("H100").Value = 45848.6
$ws.Range("J100").Value = 45848.6
$ws.Range("L100").Value = 45848.6
$ws.Range("N100").Value = -48012.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5582.8276
$ws.Range("I16").Value = 5510.15
$ws.Range("J16").Value = 5744.3335
$ws.Range("K16").Value = 5510.15
$ws.Range("L16").Value = 5744.3335
$ws.Range("M16").Value = -5223.15
$ws.Range("N16").Value = -6318.3335

$ws.Range("H107").Value = 813.0357
$ws.Range("I107").Value = 721.5
$ws.Range("K107").Value = 721.5
$ws.Range("M107").Value = 1198.5

$ws.Range("H113").Value = 5582.8276
$ws.Range("I113").Value = 5510.15
$ws.Range("J113").Value = 5744.3335
$ws.Range("K113").Value = 5510.15
$ws.Range("L113").Value = 5744.3335
$ws.Range("M113").Value = -3340.15
$ws.Range("N113").Value = -10084.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 3491.6667
$ws.Range("I115").Value = 2737.5
$ws.Range("K115").Value = 8212.5
$ws.Range("M115").Value = -7037.5

$ws.Range("H120").Value = 37500
$ws.Range("I120").Value = 12500
$ws.Range("K120").Value = 37500
$ws.Range("M120").Value = -32662

$ws.Range("H137").Value = 4111.5713
$ws.Range("I137").Value = 4465.7334
$ws.Range("J137").Value = 3226.1667
$ws.Range("K137").Value = 13397.2002
$ws.Range("L137").Value = 9678.500100000001
$ws.Range("M137").Value = -8297.200199999999
$ws.Range("N137").Value = -19878.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 53.625
$ws.Range("I2").Value = 52
$ws.Range("K2").Value = 52
$ws.Range("M2").Value = 61

$ws.Range("H132").Value = 4364.9585
$ws.Range("I132").Value = 4010.6875
$ws.Range("J132").Value = 5073.5
$ws.Range("K132").Value = 12032.0625
$ws.Range("L132").Value = 15220.5
$ws.Range("M132").Value = -9502.0625
$ws.Range("N132").Value = -20280.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5393.25
$ws.Range("I40").Value = 3998.75
$ws.Range("K40").Value = 3998.75
$ws.Range("M40").Value = -3862.75

$ws.Range("H55").Value = 2148.75
$ws.Range("I55").Value = 792.25
$ws.Range("K55").Value = 792.25
$ws.Range("M55").Value = -619.25

$ws.Range("H68").Value = 2870
$ws.Range("I68").Value = 2355.5557
$ws.Range("K68").Value = 2355.5557
$ws.Range("M68").Value = -1606.5557

$ws.Range("H71").Value = 2870
$ws.Range("I71").Value = 2355.5557
$ws.Range("K71").Value = 11777.7785
$ws.Range("M71").Value = -8033.7785

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 4511
$ws.Range("I33").Value = 4511
$ws.Range("K33").Value = 4511
$ws.Range("M33").Value = -4261

$ws.Range("H36").Value = 4511
$ws.Range("I36").Value = 4511
$ws.Range("K36").Value = 4511
$ws.Range("M36").Value = -4261

$ws.Range("H41").Value = 12007.857
$ws.Range("J41").Value = 12007.857
$ws.Range("L41").Value = 12007.857
$ws.Range("N41").Value = -12787.857

$ws.Range("H62").Value = 9954
$ws.Range("I62").Value = 13162
$ws.Range("J62").Value = 8528.223
$ws.Range("K62").Value = 13162
$ws.Range("L62").Value = 8528.223
$ws.Range("M62").Value = -12538
$ws.Range("N62").Value = -9776.223

$ws.Range("H65").Value = 9954
$ws.Range("I65").Value = 13162
$ws.Range("J65").Value = 8528.223
$ws.Range("K65").Value = 65810
$ws.Range("L65").Value = 42641.115
$ws.Range("M65").Value = -62690
$ws.Range("N65").Value = -48881.115

$ws.Range("H122").Value = 5514.364
$ws.Range("I122").Value = 5065.8
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 15197.4
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -12747.4
$ws.Range("N122").Value = -34900
